$wb = $excel.ActiveWorkbook

# --- Sheet "Daily": update row 2 (G2:L2) ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2605.47
$daily.Range("H2").Value = 5831.68
$daily.Range("I2").Value = 680.24
$daily.Range("J2").Value = 874.98
$daily.Range("K2").Value = 101.69
$daily.Range("L2").Value = 848.16

# --- Sheet "Hourly": update rows 9-19 ---
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9
$hourly.Range("K9").Value = 1.47
$hourly.Range("M9").Value = 1.47

# Row 10
$hourly.Range("K10").Value = 41.8
$hourly.Range("L10").Value = 13.13
$hourly.Range("M10").Value = 40.33

# Row 11
$hourly.Range("K11").Value = 57.79
$hourly.Range("M11").Value = 57.79

# Row 12
$hourly.Range("I12").Value = 708.35
$hourly.Range("K12").Value = 86.63
$hourly.Range("M12").Value = 86.63

# Row 13
$hourly.Range("H13").Value = 411.54
$hourly.Range("I13").Value = 760.5700000000001
$hourly.Range("K13").Value = 107.19
$hourly.Range("L13").Value = 0
$hourly.Range("M13").Value = 107.19

# Row 14
$hourly.Range("H14").Value = 442.21
$hourly.Range("I14").Value = 778.11
$hourly.Range("J14").Value = 92.55
$hourly.Range("K14").Value = 121.03
$hourly.Range("L14").Value = 0
$hourly.Range("M14").Value = 121.03

# Row 15
$hourly.Range("H15").Value = 420.94
$hourly.Range("I15").Value = 766.25
$hourly.Range("K15").Value = 125.65
$hourly.Range("L15").Value = 0
$hourly.Range("M15").Value = 125.65

# Row 16
$hourly.Range("H16").Value = 350.16
$hourly.Range("I16").Value = 721.36
$hourly.Range("J16").Value = 83.81
$hourly.Range("K16").Value = 143.96
$hourly.Range("L16").Value = 16.45
$hourly.Range("M16").Value = 138.1

# Row 17
$hourly.Range("H17").Value = 238.68
$hourly.Range("I17").Value = 627.16
$hourly.Range("J17").Value = 70.91
$hourly.Range("K17").Value = 132.55
$hourly.Range("L17").Value = 60.6
$hourly.Range("M17").Value = 115.19

# Row 18
$hourly.Range("H18").Value = 104.6
$hourly.Range("I18").Value = 427.93
$hourly.Range("K18").Value = 53.41
$hourly.Range("L18").Value = 11.52
$hourly.Range("M18").Value = 51.29

# Row 19
$hourly.Range("I19").Value = 54.9
$hourly.Range("K19").Value = 3.51
$hourly.Range("M19").Value = 3.51

$wb.Save()
